$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 75
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 6
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 12
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 5
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 4
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 10
$ws.Range("P5").Value = 0
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0

$ws.Range("W5").Select()
